$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 54. This pushes the existing row 54
# ("Vaccine Hesitancy Model") down to row 55 and existing row 55
# ("Vector Autoregression") down to row 56, along with all their
# formatting/values intact.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the "Vaccine Equity" entry.
$ws.Range("A54").Value = "Vaccine Equity"
$ws.Range("B54").Value = "A tool for exploring the effects of varying rates of vaccination among disparate socio-economic groups"
$ws.Range("C54").Value = "Damon Toth"
$ws.Range("D54").Value = "damon.toth@hcs.utah.edu"
$ws.Range("E54").Value = "Yes"
$ws.Range("G54").Value = "On development"
$ws.Range("H54").Value = "MIT"
$ws.Range("I54").Value = "R"
$ws.Range("J54").Value = "Modelers"
$ws.Range("K54").Value = "TBD"
$ws.Range("L54").Value = "Epidemic Model - Scenario Modeling"
$ws.Range("M54").Value = "Parameter inputs for simulating the model"
$ws.Range("N54").Value = "https://github.com/EpiForeSITE/vaccine-equity-model"
$ws.Range("O54").Value = "https://github.com/EpiForeSITE/vaccine-equity-model"
